# Apply cryptos list update (prices + 1h volume change) per commit diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D holds price text that sometimes looks like a plain number
# (e.g. "581.00", "2.91"). Writing such text via COM lets Excel
# auto-convert it to a numeric cell, which both changes the stored
# type and drops formatting like trailing zeros. Force the cell to
# Text format before the write, then restore the original (General)
# style so no stray number-format style is left behind.
function Set-TextValue($cell, $value) {
    $rng = $ws.Range($cell)
    $rng.NumberFormat = "@"
    $rng.Value = $value
    $rng.Style = "Normal"
}

Set-TextValue "D2" "62.260.77"
$ws.Range("E2").Value = "  -2.28%  "
Set-TextValue "D3" "3.000.05"
$ws.Range("E3").Value = "  -2.56%  "
$ws.Range("E4").Value = "  -0.03%  "
Set-TextValue "D5" "581.00"
$ws.Range("E5").Value = "  -1.20%  "
Set-TextValue "D6" "146.83"
$ws.Range("E6").Value = "  -5.48%  "
$ws.Range("E7").Value = "  -0.01%  "
Set-TextValue "D8" "0.522"
$ws.Range("E8").Value = "  -3.05%  "
Set-TextValue "D9" "3.002.93"
$ws.Range("E9").Value = "  -2.45%  "
$ws.Range("E10").Value = "  -5.39%  "
Set-TextValue "D11" "5.64"
$ws.Range("E11").Value = "  -3.59%  "
$ws.Range("E12").Value = "  -2.30%  "
$ws.Range("E13").Value = "  -4.32%  "
Set-TextValue "D14" "34.59"
$ws.Range("E14").Value = "  -5.71%  "
$ws.Range("E15").Value = "  +1.45%  "
Set-TextValue "D16" "3.498.31"
$ws.Range("E16").Value = "  -2.51%  "
$ws.Range("E17").Value = "  -2.17%  "
Set-TextValue "D18" "62.270.19"
$ws.Range("E18").Value = "  -2.18%  "
Set-TextValue "D19" "3.001.37"
$ws.Range("E19").Value = "  -2.57%  "
Set-TextValue "D20" "453.11"
$ws.Range("E20").Value = "  -3.41%  "
Set-TextValue "D21" "13.85"
$ws.Range("E21").Value = "  -3.10%  "
Set-TextValue "D22" "0.677"
$ws.Range("E22").Value = "  -3.88%  "
$ws.Range("E23").Value = "  -2.84%  "
Set-TextValue "D24" "2.29"
$ws.Range("E24").Value = "  -5.68%  "
Set-TextValue "D25" "79.97"
$ws.Range("E25").Value = "  -0.56%  "
Set-TextValue "D26" "12.26"
$ws.Range("E26").Value = "  -4.52%  "
Set-TextValue "D27" "10.07"
$ws.Range("E27").Value = "  -3.51%  "
Set-TextValue "D28" "0.999"
$ws.Range("E28").Value = "  -0.13%  "
$ws.Range("E29").Value = "  +0.00%  "
Set-TextValue "D30" "7.13"
$ws.Range("E30").Value = "  -4.26%  "
$ws.Range("E31").Value = "  -1.96%  "
$ws.Range("E32").Value = "  -1.92%  "
$ws.Range("E33").Value = "  -0.66%  "
$ws.Range("E34").Value = "  -5.12%  "
$ws.Range("E35").Value = "  -1.76%  "
Set-TextValue "D36" "0.0₃0790"
$ws.Range("E36").Value = "  -4.44%  "
$ws.Range("E37").Value = "  -4.25%  "
Set-TextValue "D38" "2.12"
$ws.Range("E38").Value = "  -3.57%  "
Set-TextValue "D39" "50.26"
$ws.Range("E39").Value = "  -0.56%  "
Set-TextValue "D40" "8.91"
$ws.Range("E40").Value = "  -2.31%  "
Set-TextValue "D41" "2.91"
$ws.Range("E41").Value = "  -10.11%  "
Set-TextValue "D42" "411.56"
$ws.Range("E42").Value = "  -4.82%  "
$ws.Range("E43").Value = "  -0.24%  "
Set-TextValue "D44" "0.276"
$ws.Range("E44").Value = "  -4.75%  "
$ws.Range("E45").Value = "  -2.58%  "
Set-TextValue "D46" "2.757.65"
$ws.Range("E46").Value = "  -1.94%  "
Set-TextValue "D47" "38.17"
$ws.Range("E47").Value = "  -4.30%  "
Set-TextValue "D48" "128.07"
$ws.Range("E48").Value = "  -1.44%  "
$ws.Range("E50").Value = "  -1.37%  "
Set-TextValue "D51" "23.67"
$ws.Range("E51").Value = "  -4.77%  "
